$wb = $excel.ActiveWorkbook

# The "Repayment Schedule" sheet gains a new (empty) column inserted
# before the existing "Late" column, pushing "Late" and "Outstanding"
# one column to the right (N -> O, O -> P, P -> Q).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and select cell S5 on it,
# which also removes the previous tab-selected state from "Transactions".
$ws.Activate()
$ws.Range("S5").Select()
